# Generate Report for Archive
#
# Refresh the localization-status report: two source files that were
# previously "Ready for handoff" have since moved back into translation,
# so their Status is updated to "In Translation" on the per-language
# detail sheets (zh-cn, de-de) as well as on the Overview roll-up sheet.
#   - 245a234a-64d5-4082-89e8-36023bd61bf6.md  -> In Translation
#   - 924bf55e-0f2f-4d64-8b4d-c4611ce3330d.md  -> In Translation
# bcf203c4-2940-446f-b99b-bf1ac4370108.md stays "Ready for handoff".

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: File Name | zh-cn | de-de -------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = $newStatus   # 245a234a...md (zh-cn)
$overview.Range("C8").Value = $newStatus   # 245a234a...md (de-de)
$overview.Range("B9").Value = $newStatus   # 924bf55e...md (zh-cn)
$overview.Range("C9").Value = $newStatus   # 924bf55e...md (de-de)

# --- Per-language detail sheets: Source File Name | Status | ... -------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B8").Value = $newStatus       # 245a234a...md
$zhcn.Range("B9").Value = $newStatus       # 924bf55e...md

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B8").Value = $newStatus       # 245a234a...md
$dede.Range("B9").Value = $newStatus       # 924bf55e...md
